$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet3"

$ws.Range("A2").Value = 30.0
$ws.Range("B2").Value = 1.0
$ws.Range("D2").Value = 29.455
$ws.Range("E2").Value = 29.455
$ws.Range("A3").Value = 30.0
$ws.Range("B3").Value = 2.0
$ws.Range("D3").Value = 29.5
$ws.Range("E3").Value = 29.5
$ws.Range("A4").Value = 30.0
$ws.Range("B4").Value = 3.0
$ws.Range("D4").Value = 29.5
$ws.Range("E4").Value = 29.5
$ws.Range("A5").Value = 30.0
$ws.Range("B5").Value = 5.0
$ws.Range("D5").Value = 29.5
$ws.Range("E5").Value = 29.5
$ws.Range("A6").Value = 30.0
$ws.Range("B6").Value = 8.0
$ws.Range("D6").Value = 29.5
$ws.Range("E6").Value = 29.5
$ws.Range("A7").Value = 67.0
$ws.Range("B7").Value = 0.0
$ws.Range("D7").Value = 34.89925
$ws.Range("E7").Value = 34.89925
$ws.Range("A8").Value = 84.0
$ws.Range("B8").Value = 2.0
$ws.Range("D8").Value = 37.411
$ws.Range("E8").Value = 37.411
$ws.Range("A9").Value = 92.0
$ws.Range("B9").Value = 3.0
$ws.Range("D9").Value = 38.8
$ws.Range("E9").Value = 38.8
$ws.Range("A10").Value = 109.0
$ws.Range("B10").Value = 5.0
$ws.Range("D10").Value = 41.349999999999994
$ws.Range("E10").Value = 41.349999999999994
$ws.Range("A11").Value = 120.0
$ws.Range("B11").Value = 8.0
$ws.Range("D11").Value = 43.0
$ws.Range("E11").Value = 43.0
$ws.Range("A12").Value = 130.0
$ws.Range("B12").Value = 1.0
$ws.Range("D12").Value = 44.11
$ws.Range("E12").Value = 44.11
$ws.Range("A13").Value = 148.0
$ws.Range("B13").Value = 2.0
$ws.Range("D13").Value = 46.756
$ws.Range("E13").Value = 46.756
$ws.Range("A14").Value = 175.0
$ws.Range("B14").Value = 3.0
$ws.Range("D14").Value = 50.724999999999994
$ws.Range("E14").Value = 50.724999999999994
$ws.Range("A15").Value = 179.0
$ws.Range("B15").Value = 4.0
$ws.Range("D15").Value = 51.849999999999994
$ws.Range("E15").Value = 51.849999999999994
$ws.Range("A16").Value = 128.0
$ws.Range("B16").Value = 9.0
$ws.Range("D16").Value = 44.2
$ws.Range("E16").Value = 44.2
$ws.Range("A17").Value = 203.0
$ws.Range("B17").Value = 1.0
$ws.Range("D17").Value = 54.68875
$ws.Range("E17").Value = 54.68875
$ws.Range("A18").Value = 202.0
$ws.Range("B18").Value = 2.0
$ws.Range("D18").Value = 54.5425
$ws.Range("E18").Value = 54.5425
$ws.Range("A19").Value = 221.0
$ws.Range("B19").Value = 3.0
$ws.Range("D19").Value = 57.32125
$ws.Range("E19").Value = 57.32125
$ws.Range("A20").Value = 243.0
$ws.Range("B20").Value = 4.0
$ws.Range("D20").Value = 61.449999999999996
$ws.Range("E20").Value = 61.449999999999996
$ws.Range("A21").Value = 300.0
$ws.Range("B21").Value = 10.0
$ws.Range("D21").Value = 70.0
$ws.Range("E21").Value = 70.0
$ws.Range("A22").Value = 549.0
$ws.Range("B22").Value = 1.0
$ws.Range("D22").Value = 104.8795
$ws.Range("E22").Value = 104.8795
$ws.Range("A23").Value = 582.0
$ws.Range("B23").Value = 2.0
$ws.Range("D23").Value = 109.681
$ws.Range("E23").Value = 109.681
$ws.Range("A24").Value = 1051.0
$ws.Range("B24").Value = 3.0
$ws.Range("D24").Value = 177.9205
$ws.Range("E24").Value = 177.9205
$ws.Range("A25").Value = 1178.0
$ws.Range("B25").Value = 4.0
$ws.Range("D25").Value = 196.399
$ws.Range("E25").Value = 196.399
$ws.Range("A26").Value = 1456.0
$ws.Range("B26").Value = 10.0
$ws.Range("D26").Value = 243.4
$ws.Range("E26").Value = 243.4
